$d = $word.ActiveDocument

# Locate the second "ASDFASDFASDF" paragraph; the very next paragraph
# (currently an empty <w:p/>) is the one that gets the new content.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "ASDFASDFASDF") {
        $targetIndex = $i
    }
}

$emptyPara = $d.Paragraphs.Item($targetIndex + 1)
$r = $emptyPara.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>All the arguments were – ‘positive’ e.g. in favor of the shit …. All the results were… ‘positive’ e.g., in favor of the shit – ERGO a plausible conclusion is that because the arguments ended up ‘having the final say’ perhaps the effect of social consensus is not actually ‘all that strong’ or ‘moral conviction’ (e.g., the FLAVOR or TYPING of the argument we used) isn’t that important, just that they got a clear argument in favor at all???</w:t></w:r><w:r><w:br/><w:t>Study 1 showed that in a vacuum however, social consensus does cohere to public opinion pretty well.</w:t></w:r></w:p>'

$r.InsertXML($xml)
